$wb = $excel.ActiveWorkbook

# Update the "Count" header to "Record Count" on the ValidationResults sheet
$wsValidation = $wb.Worksheets.Item("ValidationResults")
$wsValidation.Range("A3").Value = "Record Count"

# Make ValidationResults the active sheet with A4 selected
$wsValidation.Activate()
$wsValidation.Range("A4").Select()
